$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "ProgramsTab" query (cell B2) so the "Website" column is now
# derived from a CASE expression (falls back to program_link / program_acronym)
# instead of a straight prg.website reference.
$newProgramsQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
  CASE
    WHEN prg.program_link IS NOT NULL THEN prg.program_acronym
        ELSE prg.program_link
    END  AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
 CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type LIKE '%Lymphoma%'
ORDER BY 
    lower(prg.program_name) ASC
LIMIT 100;
'@

$ws.Range("B2").Value = $newProgramsQuery

# Scroll the view back to the top and move the selection to B2 (matching the
# refreshed sheet view saved after the edit).
$ws.Range("B2").Select()
